# Textbox response formatting fix
# Renames the 5 task-order sheets and refreshes the generated stim-file
# names (timestamps) referenced in column B of each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511687165129304"
$ws1.Range("B2").Value = "go_stims-1651168716473123.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687164948964.csv"
$ws1.Range("B4").Value = "go_stims-16511687164969616.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168716511974.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16511687189706054"
$ws2.Range("B2").Value = "OB-16511687172089832.csv"
$ws2.Range("B3").Value = "OB-16511687170086288.csv"
$ws2.Range("B4").Value = "ZB-match_6-1651168716940146.csv"
$ws2.Range("B5").Value = "TB-1651168718953321.csv"
$ws2.Range("B6").Value = "TB-16511687181379597.csv"
$ws2.Range("B7").Value = "ZB-match_6-1651168716639824.csv"
$ws2.Range("B8").Value = "OB-16511687171395278.csv"
$ws2.Range("B9").Value = "TB-1651168718493344.csv"
$ws2.Range("B10").Value = "ZB-match_9-1651168716957474.csv"

# --- Sheet 3: RS (only the sheet name changes) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651168718971574"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16511687190329437"
$ws4.Range("B2").Value = "MM_stims-16511687189861693.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687189735703.csv"
$ws4.Range("B4").Value = "MM_stims-16511687190166318.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687189861693.csv"
$ws4.Range("B6").Value = "MM_stims-1651168719031942.csv"
$ws4.Range("B7").Value = "ZM_stims-1651168719017602.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1651168719110129"
$ws5.Range("B2").Value = "vSAT_stims-1651168719095117.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687190369403.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687190471952.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511687190631936.csv"
